$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("C1").Value = "rules"
$ws.Range("E1").Value = "adaptive_filter"

# Replace numeric "1" flag in column E with the string "RLS", and refine
# the floating point precision of the F (RMSE), G (NDEI) and H (MAE)
# columns for every data row (2-16).

$values = @{
  2  = @(735.6847272333264, 1.5593073014016, 621.1362911226581)
  3  = @(628.0697862720583, 1.331213993264204, 519.8511024528511)
  4  = @(447.8658335352644, 0.9492659537818059, 371.5649217648865)
  5  = @(166.6390935828085, 0.3531968867966014, 134.1974478075941)
  6  = @(159.6070175758485, 0.3382921768634481, 123.1688304172969)
  7  = @(158.5896058539086, 0.3361357401891153, 122.237965446326)
  8  = @(156.9971585299137, 0.3327604971706271, 121.0620171341886)
  9  = @(154.8601492818763, 0.3282310377428651, 119.5049136835853)
  10 = @(152.4267656778785, 0.3230734033916748, 117.7652971102489)
  11 = @(150.051039096634, 0.3180379749436767, 115.8753280695907)
  12 = @(147.9940637366596, 0.313678149900586, 113.4334052777797)
  13 = @(146.5696033874137, 0.3106589606461219, 110.9288156094523)
  14 = @(145.6953231336185, 0.3088058956948608, 108.7440558095911)
  15 = @(144.6603842420427, 0.3066123096240889, 107.6406888505345)
  16 = @(144.8447666333826, 0.3070031139977878, 107.9718208199631)
}

foreach ($row in 2..16) {
  $ws.Cells.Item($row, 5).Value = "RLS"
  $v = $values[$row]
  $ws.Cells.Item($row, 6).Value = $v[0]
  $ws.Cells.Item($row, 7).Value = $v[1]
  $ws.Cells.Item($row, 8).Value = $v[2]
}
